$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 168.28572
$ws.Range("I8").Value = 168.28572
$ws.Range("K8").Value = 504.85716
$ws.Range("M8").Value = -365.85716

$ws.Range("H31").Value = 1687.3846
$ws.Range("I31").Value = 1096.1
$ws.Range("J31").Value = 3658.3333
$ws.Range("K31").Value = 3288.3
$ws.Range("L31").Value = 10974.9999
$ws.Range("M31").Value = -3058.3
$ws.Range("N31").Value = -11434.9999

$ws.Range("H125").Value = 12357.786
$ws.Range("I125").Value = 42581.332
$ws.Range("J125").Value = 4115
$ws.Range("K125").Value = 383231.988
$ws.Range("L125").Value = 37035
$ws.Range("M125").Value = -380771.988
$ws.Range("N125").Value = -41955

$ws.Range("H141").Value = 5667.5
$ws.Range("I141").Value = 2025.4348
$ws.Range("J141").Value = 22421
$ws.Range("K141").Value = 6076.3044
$ws.Range("L141").Value = 67263
$ws.Range("M141").Value = -896.3044
$ws.Range("N141").Value = -77623

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 602
$ws.Range("J101").Value = 602
$ws.Range("L101").Value = 602
$ws.Range("N101").Value = -7092

$ws.Range("H110").Value = 4003.1428
$ws.Range("I110").Value = 2625.2
$ws.Range("J110").Value = 7448
$ws.Range("K110").Value = 2625.2
$ws.Range("L110").Value = 7448
$ws.Range("M110").Value = -580.1999999999998
$ws.Range("N110").Value = -11538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1772.1
$ws.Range("I94").Value = 1130.1666
$ws.Range("K94").Value = 1130.1666
$ws.Range("M94").Value = -679.1666

$ws.Range("H105").Value = 4203
$ws.Range("I105").Value = 5171.6665
$ws.Range("K105").Value = 5171.6665
$ws.Range("M105").Value = -3424.6665

$ws.Range("H134").Value = 2489.4375
$ws.Range("I134").Value = 2000.2727
$ws.Range("J134").Value = 3565.6
$ws.Range("K134").Value = 6000.8181
$ws.Range("L134").Value = 10696.8
$ws.Range("M134").Value = -3465.8181
$ws.Range("N134").Value = -15766.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 72621.42999999999
$ws.Range("I4").Value = 250075
$ws.Range("J4").Value = 1640
$ws.Range("K4").Value = 250075
$ws.Range("L4").Value = 1640
$ws.Range("M4").Value = -249963
$ws.Range("N4").Value = -1864

$ws.Range("H58").Value = 1844.0741
$ws.Range("I58").Value = 823.5789
$ws.Range("J58").Value = 4267.75
$ws.Range("K58").Value = 823.5789
$ws.Range("L58").Value = 4267.75
$ws.Range("M58").Value = -620.5789
$ws.Range("N58").Value = -4673.75

$ws.Range("H136").Value = 1844.0741
$ws.Range("I136").Value = 823.5789
$ws.Range("J136").Value = 4267.75
$ws.Range("K136").Value = 2470.7367
$ws.Range("L136").Value = 12803.25
$ws.Range("M136").Value = 79.26330000000007
$ws.Range("N136").Value = -17903.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1073.0227
$ws.Range("I4").Value = 177.28572
$ws.Range("J4").Value = 1242.4865
$ws.Range("K4").Value = 531.85716
$ws.Range("L4").Value = 3727.4595
$ws.Range("M4").Value = -419.85716
$ws.Range("N4").Value = -3951.4595

$ws.Range("H11").Value = 87.0625
$ws.Range("I11").Value = 87.0625
$ws.Range("K11").Value = 261.1875
$ws.Range("M11").Value = -121.1875

$ws.Range("H98").Value = 519.13336
$ws.Range("I98").Value = 395.53845
$ws.Range("K98").Value = 1186.61535
$ws.Range("M98").Value = 311.38465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 107
$ws.Range("N5").ClearContents()

$ws.Range("H80").Value = 3093.75
$ws.Range("I80").Value = 2374.2856
$ws.Range("J80").Value = 4101
$ws.Range("K80").Value = 2374.2856
$ws.Range("L80").Value = 4101
$ws.Range("M80").Value = -1376.2856
$ws.Range("N80").Value = -6097

$ws.Range("H83").Value = 3093.75
$ws.Range("I83").Value = 2374.2856
$ws.Range("J83").Value = 4101
$ws.Range("K83").Value = 11871.428
$ws.Range("L83").Value = 20505
$ws.Range("M83").Value = -6879.428
$ws.Range("N83").Value = -30489

$ws.Range("H97").Value = 535.6667
$ws.Range("I97").Value = 474.09525
$ws.Range("K97").Value = 474.09525
$ws.Range("M97").Value = 21.90474999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4142.8213
$ws.Range("I2").Value = 999
$ws.Range("J2").Value = 4259.2593
$ws.Range("K2").Value = 999
$ws.Range("L2").Value = 4259.2593
$ws.Range("N2").Value = -4483.2593
$ws.Range("M2").Value = -887

$ws.Range("H82").Value = 1494.4
$ws.Range("I82").Value = 1491
$ws.Range("J82").Value = 1496.6666
$ws.Range("K82").Value = 1491
$ws.Range("L82").Value = 1496.6666
$ws.Range("M82").Value = -1130
$ws.Range("N82").Value = -2218.6666

$ws.Range("H85").Value = 1494.4
$ws.Range("I85").Value = 1491
$ws.Range("J85").Value = 1496.6666
$ws.Range("K85").Value = 1491
$ws.Range("L85").Value = 1496.6666
$ws.Range("M85").Value = -243
$ws.Range("N85").Value = -3992.6666

$ws.Range("H100").Value = 1253.5883
$ws.Range("I100").Value = 1114.0667
$ws.Range("J100").Value = 2300
$ws.Range("K100").Value = 1114.0667
$ws.Range("L100").Value = 2300
$ws.Range("M100").Value = -573.0667000000001
$ws.Range("N100").Value = -3382

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H100").Value = 416.66666
$ws.Range("I100").Value = 266.66666
$ws.Range("K100").Value = 533.33332
$ws.Range("M100").Value = 7.666680000000042

$ws.Range("H109").Value = 24188.5
$ws.Range("J109").Value = 24188.5
$ws.Range("L109").Value = 24188.5
$ws.Range("N109").Value = -26962.5

$ws.Range("H121").Value = 37420
$ws.Range("J121").Value = 37420
$ws.Range("L121").Value = 37420
$ws.Range("N121").Value = -40914
